$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2199312714776632
$ws.Range("C2").Value = 0.4948453608247423
$ws.Range("J2").Value = 0.01718213058419244
$ws.Range("P2").Value = 0.1683848797250859
$ws.Range("S2").Value = 0.09965635738831616
$ws.Range("B3").Value = 0.01379310344827586
$ws.Range("C3").Value = 0.01379310344827586
$ws.Range("J3").Value = 0.02758620689655172
$ws.Range("P3").Value = 0.7655172413793103
$ws.Range("S3").Value = 0.1793103448275862
$ws.Range("J4").Value = 0.06896551724137931
$ws.Range("P4").Value = 0.6896551724137931
$ws.Range("S4").Value = 0.2413793103448276
$ws.Range("B6").Value = 0.07391304347826087
$ws.Range("D6").Value = 0.01304347826086956
$ws.Range("F6").Value = 0.05217391304347826
$ws.Range("J6").Value = 0.2434782608695652
$ws.Range("O6").Value = 0.03478260869565217
$ws.Range("Q6").Value = 0.1217391304347826
$ws.Range("R6").Value = 0.08260869565217391
$ws.Range("S6").Value = 0.3782608695652174
$ws.Range("B7").Value = 0.1173708920187793
$ws.Range("D7").Value = 0.004694835680751174
$ws.Range("F7").Value = 0.06103286384976526
$ws.Range("J7").Value = 0.07981220657276995
$ws.Range("O7").Value = 0.02816901408450704
$ws.Range("Q7").Value = 0.1690140845070423
$ws.Range("R7").Value = 0.08450704225352113
$ws.Range("S7").Value = 0.4553990610328639
$ws.Range("B8").Value = 0.08528784648187633
$ws.Range("D8").Value = 0.004264392324093817
$ws.Range("E8").Value = 0.002132196162046908
$ws.Range("F8").Value = 0.06823027718550106
$ws.Range("J8").Value = 0.1044776119402985
$ws.Range("O8").Value = 0.01918976545842218
$ws.Range("Q8").Value = 0.1492537313432836
$ws.Range("R8").Value = 0.09381663113006397
$ws.Range("S8").Value = 0.4733475479744136
$ws.Range("B9").Value = 0.06862745098039216
$ws.Range("F9").Value = 0.06862745098039216
$ws.Range("J9").Value = 0.08823529411764706
$ws.Range("O9").Value = 0.04901960784313725
$ws.Range("Q9").Value = 0.1470588235294118
$ws.Range("R9").Value = 0.1666666666666667
$ws.Range("S9").Value = 0.4117647058823529
$ws.Range("B10").Value = 0.1152815013404826
$ws.Range("D10").Value = 0.02234137622877569
$ws.Range("E10").Value = 0.0008936550491510277
$ws.Range("F10").Value = 0.07059874888293119
$ws.Range("J10").Value = 0.1018766756032172
$ws.Range("O10").Value = 0.01072386058981233
$ws.Range("Q10").Value = 0.194816800714924
$ws.Range("R10").Value = 0.08936550491510277
$ws.Range("S10").Value = 0.3941018766756032
$ws.Range("G11").Value = 0.1810089020771513
$ws.Range("J11").Value = 0.09792284866468842
$ws.Range("K11").Value = 0.2522255192878338
$ws.Range("L11").Value = 0.4540059347181009
$ws.Range("S11").Value = 0.01483679525222552
$ws.Range("G12").Value = 0.78125
$ws.Range("J12").Value = 0.15
$ws.Range("K12").Value = 0.01875
$ws.Range("L12").Value = 0.025
$ws.Range("S12").Value = 0.025
$ws.Range("G13").Value = 0.8333333333333334
$ws.Range("J13").Value = 0.1428571428571428
$ws.Range("S13").Value = 0.02380952380952381
$ws.Range("F15").Value = 0.02479338842975207
$ws.Range("H15").Value = 0.2024793388429752
$ws.Range("I15").Value = 0.04132231404958678
$ws.Range("J15").Value = 0.371900826446281
$ws.Range("K15").Value = 0.1074380165289256
$ws.Range("M15").Value = 0.008264462809917356
$ws.Range("N15").Value = 0.004132231404958678
$ws.Range("O15").Value = 0.08677685950413223
$ws.Range("S15").Value = 0.1528925619834711
$ws.Range("F16").Value = 0.02298850574712644
$ws.Range("H16").Value = 0.2183908045977012
$ws.Range("I16").Value = 0.01149425287356322
$ws.Range("J16").Value = 0.4252873563218391
$ws.Range("K16").Value = 0.1264367816091954
$ws.Range("M16").Value = 0.01149425287356322
$ws.Range("O16").Value = 0.02873563218390805
$ws.Range("S16").Value = 0.1551724137931035
$ws.Range("F17").Value = 0.03532608695652174
$ws.Range("H17").Value = 0.1467391304347826
$ws.Range("I17").Value = 0.05978260869565218
$ws.Range("J17").Value = 0.3967391304347826
$ws.Range("K17").Value = 0.1005434782608696
$ws.Range("M17").Value = 0.03804347826086957
$ws.Range("O17").Value = 0.08967391304347826
$ws.Range("S17").Value = 0.1331521739130435
$ws.Range("F18").Value = 0.03015075376884422
$ws.Range("H18").Value = 0.1959798994974874
$ws.Range("I18").Value = 0.07035175879396985
$ws.Range("J18").Value = 0.3819095477386935
$ws.Range("K18").Value = 0.09045226130653267
$ws.Range("M18").Value = 0.01005025125628141
$ws.Range("O18").Value = 0.1155778894472362
$ws.Range("S18").Value = 0.1055276381909548
$ws.Range("F19").Value = 0.0290088638195004
$ws.Range("H19").Value = 0.2344883158742949
$ws.Range("I19").Value = 0.04593070104754231
$ws.Range("J19").Value = 0.346494762288477
$ws.Range("K19").Value = 0.1144238517324738
$ws.Range("M19").Value = 0.0185334407735697
$ws.Range("O19").Value = 0.0781627719580983
$ws.Range("S19").Value = 0.1329572925060435
